$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 6666799
$ws.Range("I5").Value = 6666799
$ws.Range("K5").Value = 6666799
$ws.Range("M5").Value = -6666684
$ws.Range("H53").Value = 278.45456
$ws.Range("I53").Value = 333.9375
$ws.Range("K53").Value = 333.9375
$ws.Range("M53").Value = 303.0625
$ws.Range("H74").Value = 3883.1667
$ws.Range("I74").Value = 3859.8
$ws.Range("K74").Value = 3859.8
$ws.Range("M74").Value = -2923.8
$ws.Range("H77").Value = 3883.1667
$ws.Range("I77").Value = 3859.8
$ws.Range("K77").Value = 19299
$ws.Range("M77").Value = -14619
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H87").Value = 23999.666
$ws.Range("J87").Value = 23999.666
$ws.Range("L87").Value = 23999.666
$ws.Range("N87").Value = -26495.666
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H90").Value = 23999.666
$ws.Range("J90").Value = 23999.666
$ws.Range("L90").Value = 71998.99800000001
$ws.Range("N90").Value = -84478.99800000001
$ws.Range("H112").Value = 3608.75
$ws.Range("I112").Value = 3624
$ws.Range("K112").Value = 10872
$ws.Range("M112").Value = -9764
$ws.Range("H113").Value = 7964.6924
$ws.Range("J113").Value = 9103.700000000001
$ws.Range("L113").Value = 9103.700000000001
$ws.Range("N113").Value = -15611.7
$ws.Range("H116").Value = 5259.2856
$ws.Range("J116").Value = 5347
$ws.Range("L116").Value = 5347
$ws.Range("N116").Value = -12231
$ws.Range("H137").Value = 2122.9092
$ws.Range("J137").Value = 2397.6667
$ws.Range("L137").Value = 7193.000100000001
$ws.Range("N137").Value = -12293.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1959.5333
$ws.Range("I74").Value = 1791.8
$ws.Range("K74").Value = 1791.8
$ws.Range("M74").Value = -917.8
$ws.Range("H77").Value = 1959.5333
$ws.Range("I77").Value = 1791.8
$ws.Range("K77").Value = 8959
$ws.Range("M77").Value = -4591
$ws.Range("H122").Value = 2559.0833
$ws.Range("I122").Value = 2559.0833
$ws.Range("K122").Value = 7677.249899999999
$ws.Range("M122").Value = -5227.249899999999
$ws.Range("H132").Value = 4337.154
$ws.Range("I132").Value = 3307.6365
$ws.Range("K132").Value = 9922.9095
$ws.Range("M132").Value = -7392.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3290.2144
$ws.Range("I31").Value = 2461
$ws.Range("K31").Value = 2461
$ws.Range("M31").Value = -2166
$ws.Range("H34").Value = 3290.2144
$ws.Range("I34").Value = 2461
$ws.Range("K34").Value = 2461
$ws.Range("M34").Value = -2259
$ws.Range("H54").Value = 42250
$ws.Range("J54").Value = 65000
$ws.Range("L54").Value = 65000
$ws.Range("N54").Value = -66316
$ws.Range("H58").Value = 2656.9211
$ws.Range("I58").Value = 2170.44
$ws.Range("J58").Value = 3592.4614
$ws.Range("K58").Value = 2170.44
$ws.Range("L58").Value = 3592.4614
$ws.Range("M58").Value = -1967.44
$ws.Range("N58").Value = -3998.4614
$ws.Range("H99").Value = 19463.826
$ws.Range("I99").Value = 17697.5
$ws.Range("K99").Value = 17697.5
$ws.Range("M99").Value = -16199.5
$ws.Range("H122").Value = 4548.625
$ws.Range("I122").Value = 1963.6666
$ws.Range("K122").Value = 5890.9998
$ws.Range("M122").Value = -3440.9998
$ws.Range("H126").Value = 19463.826
$ws.Range("I126").Value = 17697.5
$ws.Range("K126").Value = 53092.5
$ws.Range("M126").Value = -50622.5
$ws.Range("H136").Value = 2656.9211
$ws.Range("I136").Value = 2170.44
$ws.Range("J136").Value = 3592.4614
$ws.Range("K136").Value = 6511.32
$ws.Range("L136").Value = 10777.3842
$ws.Range("M136").Value = -3961.32
$ws.Range("N136").Value = -15877.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2823.652
$ws.Range("J34").Value = 4388.9287
$ws.Range("L34").Value = 13166.7861
$ws.Range("N34").Value = -13334.7861
$ws.Range("H39").Value = 3347.5
$ws.Range("J39").Value = 3347.5
$ws.Range("L39").Value = 10042.5
$ws.Range("N39").Value = -10630.5
$ws.Range("H113").Value = 1416.8
$ws.Range("I113").Value = 2129.8
$ws.Range("J113").Value = 1179.1333
$ws.Range("K113").Value = 6389.400000000001
$ws.Range("L113").Value = 3537.3999
$ws.Range("M113").Value = -4219.400000000001
$ws.Range("N113").Value = -7877.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 12000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H64").Value = 70271
$ws.Range("J64").Value = 70271
$ws.Range("L64").Value = 70271
$ws.Range("N64").Value = -70767
$ws.Range("H67").Value = 70271
$ws.Range("J67").Value = 70271
$ws.Range("L67").Value = 70271
$ws.Range("N67").Value = -71987
$ws.Range("H113").Value = 3888.7144
$ws.Range("I113").Value = 3706.2727
$ws.Range("J113").Value = 4557.6665
$ws.Range("K113").Value = 3706.2727
$ws.Range("L113").Value = 4557.6665
$ws.Range("M113").Value = -1536.2727
$ws.Range("N113").Value = -8897.666499999999
$ws.Range("H122").Value = 1998.2222
$ws.Range("I122").Value = 1935.5714
$ws.Range("K122").Value = 5806.7142
$ws.Range("M122").Value = -3356.7142
$ws.Range("H126").Value = 5725.5
$ws.Range("I126").Value = 5485
$ws.Range("K126").Value = 16455
$ws.Range("M126").Value = -13985
$ws.Range("H132").Value = 4993.6665
$ws.Range("I132").Value = 4993.6665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14980.9995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12450.9995
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 5518423
$ws.Range("I43").Value = 2795555.5
$ws.Range("K43").Value = 2795555.5
$ws.Range("M43").Value = -2795362.5
$ws.Range("H100").Value = 3902.875
$ws.Range("I100").Value = 2742
$ws.Range("K100").Value = 2742
$ws.Range("M100").Value = -2201
$ws.Range("H132").Value = 5300
$ws.Range("I132").Value = 3700
$ws.Range("J132").Value = 8500
$ws.Range("K132").Value = 11100
$ws.Range("L132").Value = 25500
$ws.Range("M132").Value = -8570
$ws.Range("N132").Value = -30560
